$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.124.76'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.878.83'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5063'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3848'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09043'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.126'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.367'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").Value = '1.872.04'
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.273'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001111'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06592'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.132'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").Value = '28.156.29'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.551'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.093.94'
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.060'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.623'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.518'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06608'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02408'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2200'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.294'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.216'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6443'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6047'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.665'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.273'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.239'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.008'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.95%  '
